$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich text shared strings) ---
# A8: "Volume 29   Number  47" -> "...48"  (replace just the "47" run)
$volChars = $ws.Range("A8").Characters(21, 2)
$volChars.Text = "48"

# C9: "Report Covering the Week  11/21/2022  Through  11/27/2022"
# -> "...11/28/2022...12/4/2022"
$d1Chars = $ws.Range("C9").Characters(27, 10)
$d1Chars.Text = "11/28/2022"
$d2Chars = $ws.Range("C9").Characters(48, 10)
$d2Chars.Text = "12/4/2022"

# --- Data table updates (rows 14-30) ---
$ws.Range("C14").Value = 2
$ws.Range("E14").Value = -75
$ws.Range("F14").Value = 26
$ws.Range("G14").Value = 24
$ws.Range("H14").Value = 8.333333333333
$ws.Range("I14").Value = 393
$ws.Range("J14").Value = 448
$ws.Range("K14").Value = -12.276785714285
$ws.Range("L14").Value = -10.068649885583
$ws.Range("M14").Value = -21.556886227544
$ws.Range("N14").Value = -78.032420346562
$ws.Range("C15").Value = 26
$ws.Range("D15").Value = 30
$ws.Range("E15").Value = -13.333333333333
$ws.Range("F15").Value = 97
$ws.Range("G15").Value = 117
$ws.Range("H15").Value = -17.094017094017
$ws.Range("I15").Value = 1512
$ws.Range("J15").Value = 1395
$ws.Range("K15").Value = 8.387096774193
$ws.Range("L15").Value = 11.094783247612
$ws.Range("M15").Value = 17.757009345794
$ws.Range("N15").Value = -50.049554013875
$ws.Range("C16").Value = 298
$ws.Range("D16").Value = 366
$ws.Range("E16").Value = -18.579234972677
$ws.Range("F16").Value = 1357
$ws.Range("G16").Value = 1361
$ws.Range("H16").Value = -0.293901542983
$ws.Range("I16").Value = 16274
$ws.Range("J16").Value = 12693
$ws.Range("K16").Value = 28.212400535728
$ws.Range("L16").Value = 34.185356200527
$ws.Range("M16").Value = -10.003871039097
$ws.Range("N16").Value = -79.525696672328
$ws.Range("C17").Value = 409
$ws.Range("D17").Value = 443
$ws.Range("E17").Value = -7.674943566591
$ws.Range("F17").Value = 1784
$ws.Range("G17").Value = 1821
$ws.Range("H17").Value = -2.031850631521
$ws.Range("I17").Value = 24208
$ws.Range("J17").Value = 21429
$ws.Range("K17").Value = 12.968407298520
$ws.Range("L17").Value = 24.815674142820
$ws.Range("M17").Value = 51.546262676849
$ws.Range("N17").Value = -37.479338842975
$ws.Range("C18").Value = 301
$ws.Range("D18").Value = 323
$ws.Range("E18").Value = -6.811145510835
$ws.Range("F18").Value = 1165
$ws.Range("G18").Value = 1238
$ws.Range("H18").Value = -5.896607431340
$ws.Range("I18").Value = 14602
$ws.Range("J18").Value = 11674
$ws.Range("K18").Value = 25.081377419907
$ws.Range("L18").Value = 1.269158748873
$ws.Range("M18").Value = -16.123843988741
$ws.Range("N18").Value = -84.384557801304
$ws.Range("C19").Value = 986
$ws.Range("D19").Value = 1400
$ws.Range("E19").Value = -29.571428571428
$ws.Range("F19").Value = 3748
$ws.Range("G19").Value = 4379
$ws.Range("H19").Value = -14.409682575930
$ws.Range("I19").Value = 47821
$ws.Range("J19").Value = 36326
$ws.Range("K19").Value = 31.644001541595
$ws.Range("L19").Value = 43.658375390531
$ws.Range("M19").Value = 35.839677309396
$ws.Range("N19").Value = -39.893917874335
$ws.Range("C20").Value = 287
$ws.Range("D20").Value = 210
$ws.Range("E20").Value = 36.666666666666
$ws.Range("F20").Value = 1064
$ws.Range("G20").Value = 960
$ws.Range("H20").Value = 10.833333333333
$ws.Range("I20").Value = 12596
$ws.Range("J20").Value = 9549
$ws.Range("K20").Value = 31.909100429364
$ws.Range("L20").Value = 50.256471430275
$ws.Range("M20").Value = 31.249348754819
$ws.Range("N20").Value = -87.901955511160
$ws.Range("C21").Value = 2309
$ws.Range("D21").Value = 2780
$ws.Range("E21").Value = -16.942446043165
$ws.Range("F21").Value = 9241
$ws.Range("G21").Value = 9900
$ws.Range("H21").Value = -6.656565656565
$ws.Range("I21").Value = 117406
$ws.Range("J21").Value = 93514
$ws.Range("K21").Value = 25.549115640438
$ws.Range("L21").Value = 31.310465155294
$ws.Range("M21").Value = 19.738506098804
$ws.Range("N21").Value = -70.663754847479
$ws.Range("C22").Value = 49
$ws.Range("D22").Value = 58
$ws.Range("E22").Value = -15.517241379310
$ws.Range("F22").Value = 162
$ws.Range("G22").Value = 218
$ws.Range("H22").Value = -25.688073394495
$ws.Range("I22").Value = 2147
$ws.Range("J22").Value = 1629
$ws.Range("K22").Value = 31.798649478207
$ws.Range("L22").Value = 29.728096676737
$ws.Range("M22").Value = 8.324924318869
$ws.Range("C23").Value = 95
$ws.Range("D23").Value = 128
$ws.Range("E23").Value = -25.78125
$ws.Range("F23").Value = 434
$ws.Range("G23").Value = 454
$ws.Range("H23").Value = -4.405286343612
$ws.Range("I23").Value = 5519
$ws.Range("J23").Value = 5137
$ws.Range("K23").Value = 7.436246836675
$ws.Range("L23").Value = 15.291414246918
$ws.Range("M23").Value = 40.218495934959
$ws.Range("C24").Value = 2404
$ws.Range("D24").Value = 2181
$ws.Range("E24").Value = 10.224667583677
$ws.Range("F24").Value = 8829
$ws.Range("G24").Value = 8155
$ws.Range("H24").Value = 8.264868179031
$ws.Range("I24").Value = 107193
$ws.Range("J24").Value = 79365
$ws.Range("K24").Value = 35.063315063315
$ws.Range("L24").Value = 41.60986049461
$ws.Range("M24").Value = 40.785930994628
$ws.Range("C25").Value = 679
$ws.Range("D25").Value = 712
$ws.Range("E25").Value = -4.634831460674
$ws.Range("F25").Value = 2950
$ws.Range("G25").Value = 3043
$ws.Range("H25").Value = -3.056194544857
$ws.Range("I25").Value = 38411
$ws.Range("J25").Value = 33617
$ws.Range("K25").Value = 14.260641937115
$ws.Range("L25").Value = 23.70692431562
$ws.Range("M25").Value = -10.359393232205
$ws.Range("C26").Value = 37
$ws.Range("D26").Value = 56
$ws.Range("E26").Value = -33.928571428571
$ws.Range("F26").Value = 160
$ws.Range("H26").Value = -21.951219512195
$ws.Range("I26").Value = 2422
$ws.Range("J26").Value = 2301
$ws.Range("K26").Value = 5.258583224684
$ws.Range("L26").Value = 14.299197734780
$ws.Range("C27").Value = 80
$ws.Range("D27").Value = 107
$ws.Range("E27").Value = -25.233644859813
$ws.Range("F27").Value = 363
$ws.Range("G27").Value = 407
$ws.Range("H27").Value = -10.810810810810
$ws.Range("I27").Value = 4829
$ws.Range("J27").Value = 4578
$ws.Range("K27").Value = 5.482743556138
$ws.Range("L27").Value = 35.341928251121
$ws.Range("C28").Value = 13
$ws.Range("D28").Value = 36
$ws.Range("E28").Value = -63.888888888888
$ws.Range("F28").Value = 82
$ws.Range("G28").Value = 147
$ws.Range("H28").Value = -44.217687074829
$ws.Range("I28").Value = 1474
$ws.Range("J28").Value = 1757
$ws.Range("K28").Value = -16.107000569152
$ws.Range("L28").Value = -15.723270440251
$ws.Range("M28").Value = -11.524609843937
$ws.Range("N28").Value = -73.369467028003
$ws.Range("C29").Value = 13
$ws.Range("E29").Value = -59.375
$ws.Range("F29").Value = 72
$ws.Range("G29").Value = 125
$ws.Range("H29").Value = -42.4
$ws.Range("I29").Value = 1218
$ws.Range("J29").Value = 1467
$ws.Range("K29").Value = -16.973415132924
$ws.Range("L29").Value = -14.586255259467
$ws.Range("M29").Value = -11.418181818181
$ws.Range("N29").Value = -75.488025759710
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 11
$ws.Range("E30").Value = -90.909090909090
$ws.Range("F30").Value = 26
$ws.Range("H30").Value = -21.212121212121
$ws.Range("I30").Value = 571
$ws.Range("J30").Value = 495
$ws.Range("K30").Value = 15.353535353535
$ws.Range("L30").Value = 130.241935483871
